# Apply the "Add files via upload" edit:
#   Slide 4 ("Logistic Regression Model" -> "Multiple Linear Regression Model"):
#     - give the title placeholder an explicit position/size (xfrm)
#     - split the title text into three runs: "Multiple " / "Linear " / "Regression Model"
#   Slide 5 ("Importance of Research"):
#     - merge the two title runs ("Importance " + "of Research") into a single run

$p = $ppt.ActivePresentation

# --- Slide 4: title placeholder ---
$s4  = $p.Slides.Item(4)
$sh4 = $s4.Shapes.Item(1)

# Explicit xfrm: off x=457199 y=152718, ext cx=7347593 cy=1371600 (EMU -> points, 12700 EMU/pt)
$sh4.Left   = 35.99993
$sh4.Top    = 12.02503937007874
$sh4.Width  = 578.5506299212599
$sh4.Height = 108.0

$tr4 = $sh4.TextFrame.TextRange
$tr4.Text = "Multiple "
$tr4.InsertAfter("Linear ") | Out-Null
$tr4.InsertAfter("Regression Model") | Out-Null

# --- Slide 5: title placeholder ---
$s5  = $p.Slides.Item(5)
$sh5 = $s5.Shapes.Item(1)
$tr5 = $sh5.TextFrame.TextRange

# Force a real text rewrite through an unrelated placeholder string first
# (the setter keeps any shared prefix/suffix as separate runs, so going via
# a string with nothing in common collapses the whole paragraph down to one
# run before we write the final merged text).
$tr5.Text = "ZZZZZZZZZZZZZZZZZZZZZZZ"
$tr5b = $sh5.TextFrame.TextRange
$tr5b.Text = "Importance of Research"
